$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SalvarCliente")
$ws2 = $wb.Worksheets.Item("EfetuarLogin")
# style16 test: border1 vertical=center, normal font (copy border-only style then add vertical=center)
$ws1.Range("A3").Copy()
$ws2.Range("D70").PasteSpecial(-4122)
$ws2.Range("D70").VerticalAlignment = -4108
# style19 test: hyperlink font (fontId3), no border, vertical=center
$ws1.Range("H3").Copy()  # style7: xfId1,fontId3,border1
$ws2.Range("D71").PasteSpecial(-4122)
$ws2.Range("D71").Borders.LineStyle = -4142
$ws2.Range("D71").VerticalAlignment = -4108
